$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All existing rows (2..468) got their "Förändrad" (column C) date bumped
# from 2023-09-20 (45189) to 2023-09-21 (45190).
$ws.Range("C2:C468").Value = 45190

# Row 468 (previously the last data row) picks up an explicit row height,
# matching the other data rows.
$ws.Rows.Item(468).RowHeight = 15

# A new record is appended as row 469.
$r = 469
$ws.Cells.Item($r, 1).Value = "A 44034-2023"

$ws.Cells.Item($r, 2).Value = 45188
$ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($r, 3).Value = 45190
$ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($r, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item($r, 5).Value = "SMEDJEBACKEN"
$ws.Cells.Item($r, 6).Value = "Sveaskog"
$ws.Cells.Item($r, 7).Value = 0.4
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0

# Column R keeps the wrap-text style used throughout the rest of the sheet,
# even though it has no content for this row.
$ws.Cells.Item($r, 18).WrapText = $true
